$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.090240666666667
$ws.Range("H2").Value = 15.270722
$ws.Range("I2").Value = 0.01518526656315472
$ws.Range("J2").Value = 0.01525191836740238
$ws.Range("M2").Value = 2.914938
$ws.Range("N2").Value = 8.744814
$ws.Range("O2").Value = 0.07105539873786189
$ws.Range("P2").Value = 0.07541359350096061
$ws.Range("Q2").Value = 14.837735948412
$ws.Range("R2").Value = 133.539623535708
$ws.Range("S2").Value = 0.00107899517058568
$ws.Range("T2").Value = 0.001150201971869118
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.090240666666667
$ws.Range("H3").Value = 15.270722
$ws.Range("I3").Value = 0.01518526656315472
$ws.Range("J3").Value = 0.01525191836740238
$ws.Range("M3").Value = 7.339638666666666
$ws.Range("O3").Value = 0.1789132228719201
$ws.Range("P3").Value = 0.1898868953137022
$ws.Range("Q3").Value = 37.36052721970578
$ws.Range("R3").Value = 336.244744977352
$ws.Range("S3").Value = 0.002716844980983216
$ws.Range("T3").Value = 0.002896139426364068
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.090240666666667
$ws.Range("H4").Value = 15.270722
$ws.Range("I4").Value = 0.01518526656315472
$ws.Range("J4").Value = 0.01525191836740238
$ws.Range("M4").Value = 15.66105466666667
$ws.Range("N4").Value = 46.983164
$ws.Range("O4").Value = 0.3817585430617917
$ws.Range("P4").Value = 0.4051737671361525
$ws.Range("Q4").Value = 79.71853734715646
$ws.Range("R4").Value = 717.4668361244081
$ws.Range("S4").Value = 0.005797105239154886
$ws.Range("T4").Value = 0.006179677220973501
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.090240666666667
$ws.Range("H5").Value = 15.270722
$ws.Range("I5").Value = 0.01518526656315472
$ws.Range("J5").Value = 0.01525191836740238
$ws.Range("M5").Value = 7.112307
$ws.Range("N5").Value = 14.224614
$ws.Range("O5").Value = 0.1733717183113625
$ws.Range("P5").Value = 0.1226703344295343
$ws.Range("Q5").Value = 36.20335432521801
$ws.Range("R5").Value = 217.220125951308
$ws.Range("S5").Value = 0.002632695757070212
$ws.Range("T5").Value = 0.001870957926821207
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.090240666666667
$ws.Range("H6").Value = 15.270722
$ws.Range("I6").Value = 0.01518526656315472
$ws.Range("J6").Value = 0.01525191836740238
$ws.Range("M6").Value = 7.995517333333333
$ws.Range("N6").Value = 23.986552
$ws.Range("O6").Value = 0.1949011170170639
$ws.Range("P6").Value = 0.2068554096196504
$ws.Range("Q6").Value = 40.69910748117156
$ws.Range("R6").Value = 366.291967330544
$ws.Range("S6").Value = 0.002959625415360726
$ws.Range("T6").Value = 0.003154941821374489
$ws.Range("I7").Value = 0.003043737298347591
$ws.Range("J7").Value = 0.003057096996825524
$ws.Range("M7").Value = 2.914938
$ws.Range("N7").Value = 8.744814
$ws.Range("O7").Value = 0.07105539873786189
$ws.Range("P7").Value = 0.07541359350096061
$ws.Range("Q7").Value = 2.974078205436
$ws.Range("R7").Value = 26.766703848924
$ws.Range("S7").Value = 0.0002162739673873906
$ws.Range("T7").Value = 0.0002305466702116075
$ws.Range("I8").Value = 0.003043737298347591
$ws.Range("J8").Value = 0.003057096996825524
$ws.Range("M8").Value = 7.339638666666666
$ws.Range("O8").Value = 0.1789132228719201
$ws.Range("P8").Value = 0.1898868953137022
$ws.Range("Q8").Value = 7.488550149028443
$ws.Range("R8").Value = 67.39695134125598
$ws.Range("S8").Value = 0.0005445648496228385
$ws.Range("T8").Value = 0.0005805026574000415
$ws.Range("I9").Value = 0.003043737298347591
$ws.Range("J9").Value = 0.003057096996825524
$ws.Range("M9").Value = 15.66105466666667
$ws.Range("N9").Value = 46.983164
$ws.Range("O9").Value = 0.3817585430617917
$ws.Range("P9").Value = 0.4051737671361525
$ws.Range("Q9").Value = 15.97879658444711
$ws.Range("R9").Value = 143.809169260024
$ws.Range("S9").Value = 0.00116197271648001
$ws.Range("T9").Value = 0.001238655506704416
$ws.Range("I10").Value = 0.003043737298347591
$ws.Range("J10").Value = 0.003057096996825524
$ws.Range("M10").Value = 7.112307
$ws.Range("N10").Value = 14.224614
$ws.Range("O10").Value = 0.1733717183113625
$ws.Range("P10").Value = 0.1226703344295343
$ws.Range("Q10").Value = 7.256606225954
$ws.Range("R10").Value = 43.539637355724
$ws.Range("S10").Value = 0.0005276979655029061
$ws.Range("T10").Value = 0.0003750151109841118
$ws.Range("I11").Value = 0.003043737298347591
$ws.Range("J11").Value = 0.003057096996825524
$ws.Range("M11").Value = 7.995517333333333
$ws.Range("N11").Value = 23.986552
$ws.Range("O11").Value = 0.1949011170170639
$ws.Range("P11").Value = 0.2068554096196504
$ws.Range("Q11").Value = 8.157735719336888
$ws.Range("R11").Value = 73.419621474032
$ws.Range("S11").Value = 0.0005932277993544459
$ws.Range("T11").Value = 0.0006323770515253468
$ws.Range("G12").Value = 69.13821133333333
$ws.Range("H12").Value = 207.414634
$ws.Range("I12").Value = 0.206253935235621
$ws.Range("J12").Value = 0.2071592335956769
$ws.Range("M12").Value = 2.914938
$ws.Range("N12").Value = 8.744814
$ws.Range("O12").Value = 0.07105539873786189
$ws.Range("P12").Value = 0.07541359350096061
$ws.Range("Q12").Value = 201.533599467564
$ws.Range("R12").Value = 1813.802395208076
$ws.Range("S12").Value = 0.0146554556094202
$ws.Range("T12").Value = 0.01562262223235492
$ws.Range("G13").Value = 69.13821133333333
$ws.Range("H13").Value = 207.414634
$ws.Range("I13").Value = 0.206253935235621
$ws.Range("J13").Value = 0.2071592335956769
$ws.Range("M13").Value = 7.339638666666666
$ws.Range("O13").Value = 0.1789132228719201
$ws.Range("P13").Value = 0.1898868953137022
$ws.Range("Q13").Value = 507.4494892463048
$ws.Range("R13").Value = 4567.045403216743
$ws.Range("S13").Value = 0.03690155628302124
$ws.Range("T13").Value = 0.03933682370304907
$ws.Range("G14").Value = 69.13821133333333
$ws.Range("H14").Value = 207.414634
$ws.Range("I14").Value = 0.206253935235621
$ws.Range("J14").Value = 0.2071592335956769
$ws.Range("M14").Value = 15.66105466666667
$ws.Range("N14").Value = 46.983164
$ws.Range("O14").Value = 0.3817585430617917
$ws.Range("P14").Value = 0.4051737671361525
$ws.Range("Q14").Value = 1082.777307246886
$ws.Range("R14").Value = 9744.995765221975
$ws.Range("S14").Value = 0.07873920181631183
$ws.Range("T14").Value = 0.08393548707299861
$ws.Range("G15").Value = 69.13821133333333
$ws.Range("H15").Value = 207.414634
$ws.Range("I15").Value = 0.206253935235621
$ws.Range("J15").Value = 0.2071592335956769
$ws.Range("M15").Value = 7.112307
$ws.Range("N15").Value = 14.224614
$ws.Range("O15").Value = 0.1733717183113625
$ws.Range("P15").Value = 0.1226703344295343
$ws.Range("Q15").Value = 491.732184433546
$ws.Range("R15").Value = 2950.393106601276
$ws.Range("S15").Value = 0.0357585991602801
$ws.Range("T15").Value = 0.02541229246534769
$ws.Range("G16").Value = 69.13821133333333
$ws.Range("H16").Value = 207.414634
$ws.Range("I16").Value = 0.206253935235621
$ws.Range("J16").Value = 0.2071592335956769
$ws.Range("M16").Value = 7.995517333333333
$ws.Range("N16").Value = 23.986552
$ws.Range("O16").Value = 0.1949011170170639
$ws.Range("P16").Value = 0.2068554096196504
$ws.Range("Q16").Value = 552.7957671113297
$ws.Range("R16").Value = 4975.161904001968
$ws.Range("S16").Value = 0.0401991223665877
$ws.Range("T16").Value = 0.04285200812192659
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 4.394653
$ws.Range("H17").Value = 8.789306
$ws.Range("I17").Value = 0.01311018115402158
$ws.Range("J17").Value = 0.008778483271329277
$ws.Range("M17").Value = 2.914938
$ws.Range("N17").Value = 8.744814
$ws.Range("O17").Value = 0.07105539873786189
$ws.Range("P17").Value = 0.07541359350096061
$ws.Range("Q17").Value = 12.810141026514
$ws.Range("R17").Value = 76.86084615908399
$ws.Range("S17").Value = 0.0009315491494246056
$ws.Range("T17").Value = 0.000662016968979009
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 4.394653
$ws.Range("H18").Value = 8.789306
$ws.Range("I18").Value = 0.01311018115402158
$ws.Range("J18").Value = 0.008778483271329277
$ws.Range("M18").Value = 7.339638666666666
$ws.Range("O18").Value = 0.1789132228719201
$ws.Range("P18").Value = 0.1898868953137022
$ws.Range("Q18").Value = 32.25516508538266
$ws.Range("R18").Value = 193.530990512296
$ws.Range("S18").Value = 0.002345584762700709
$ws.Range("T18").Value = 0.001666918933955988
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 4.394653
$ws.Range("H19").Value = 8.789306
$ws.Range("I19").Value = 0.01311018115402158
$ws.Range("J19").Value = 0.008778483271329277
$ws.Range("M19").Value = 15.66105466666667
$ws.Range("N19").Value = 46.983164
$ws.Range("O19").Value = 0.3817585430617917
$ws.Range("P19").Value = 0.4051737671361525
$ws.Range("Q19").Value = 68.82490087403066
$ws.Range("R19").Value = 412.949405244184
$ws.Range("S19").Value = 0.005004923656635437
$ws.Range("T19").Value = 0.003556811136786179
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 4.394653
$ws.Range("H20").Value = 8.789306
$ws.Range("I20").Value = 0.01311018115402158
$ws.Range("J20").Value = 0.008778483271329277
$ws.Range("M20").Value = 7.112307
$ws.Range("N20").Value = 14.224614
$ws.Range("O20").Value = 0.1733717183113625
$ws.Range("P20").Value = 0.1226703344295343
$ws.Range("Q20").Value = 31.256121294471
$ws.Range("R20").Value = 125.024485177884
$ws.Range("S20").Value = 0.002272934634045962
$ws.Range("T20").Value = 0.001076859478678034
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 4.394653
$ws.Range("H21").Value = 8.789306
$ws.Range("I21").Value = 0.01311018115402158
$ws.Range("J21").Value = 0.008778483271329277
$ws.Range("M21").Value = 7.995517333333333
$ws.Range("N21").Value = 23.986552
$ws.Range("O21").Value = 0.1949011170170639
$ws.Range("P21").Value = 0.2068554096196504
$ws.Range("Q21").Value = 35.13752423548533
$ws.Range("R21").Value = 210.825145412912
$ws.Range("S21").Value = 0.002555188951214866
$ws.Range("T21").Value = 0.001815876752930066
$ws.Range("G22").Value = 255.5657806666667
$ws.Range("H22").Value = 766.6973419999999
$ws.Range("I22").Value = 0.7624068797488551
$ws.Range("J22").Value = 0.765753267768766
$ws.Range("M22").Value = 2.914938
$ws.Range("N22").Value = 8.744814
$ws.Range("O22").Value = 0.07105539873786189
$ws.Range("P22").Value = 0.07541359350096061
$ws.Range("Q22").Value = 744.9584055649319
$ws.Range("R22").Value = 6704.625650084387
$ws.Range("S22").Value = 0.05417312484104402
$ws.Range("T22").Value = 0.05774820565754596
$ws.Range("G23").Value = 255.5657806666667
$ws.Range("H23").Value = 766.6973419999999
$ws.Range("I23").Value = 0.7624068797488551
$ws.Range("J23").Value = 0.765753267768766
$ws.Range("M23").Value = 7.339638666666666
$ws.Range("O23").Value = 0.1789132228719201
$ws.Range("P23").Value = 0.1898868953137022
$ws.Range("Q23").Value = 1875.760485657919
$ws.Range("R23").Value = 16881.84437092127
$ws.Range("S23").Value = 0.1364046719955921
$ws.Range("T23").Value = 0.145406510592933
$ws.Range("G24").Value = 255.5657806666667
$ws.Range("H24").Value = 766.6973419999999
$ws.Range("I24").Value = 0.7624068797488551
$ws.Range("J24").Value = 0.765753267768766
$ws.Range("M24").Value = 15.66105466666667
$ws.Range("N24").Value = 46.983164
$ws.Range("O24").Value = 0.3817585430617917
$ws.Range("P24").Value = 0.4051737671361525
$ws.Range("Q24").Value = 4002.42966195001
$ws.Range("R24").Value = 36021.86695755009
$ws.Range("S24").Value = 0.2910553396332096
$ws.Range("T24").Value = 0.3102631361986898
$ws.Range("G25").Value = 255.5657806666667
$ws.Range("H25").Value = 766.6973419999999
$ws.Range("I25").Value = 0.7624068797488551
$ws.Range("J25").Value = 0.765753267768766
$ws.Range("M25").Value = 7.112307
$ws.Range("N25").Value = 14.224614
$ws.Range("O25").Value = 0.1733717183113625
$ws.Range("P25").Value = 0.1226703344295343
$ws.Range("Q25").Value = 1817.662290795998
$ws.Range("R25").Value = 10905.97374477599
$ws.Range("S25").Value = 0.1321797907944633
$ws.Range("T25").Value = 0.09393520944770321
$ws.Range("G26").Value = 255.5657806666667
$ws.Range("H26").Value = 766.6973419999999
$ws.Range("I26").Value = 0.7624068797488551
$ws.Range("J26").Value = 0.765753267768766
$ws.Range("M26").Value = 7.995517333333333
$ws.Range("N26").Value = 23.986552
$ws.Range("O26").Value = 0.1949011170170639
$ws.Range("P26").Value = 0.2068554096196504
$ws.Range("Q26").Value = 2043.380629127198
$ws.Range("R26").Value = 18390.42566214478
$ws.Range("S26").Value = 0.1485939524845462
$ws.Range("T26").Value = 0.1584002058718939
